$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.891550288703
$ws.Range("B3").Value = 0.8047241718875351
$ws.Range("B4").Value = 0.08962212772212796
$ws.Range("B5").Value = 0.2934480393947849
